$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 257, shifting existing rows 257:314 down to 258:315.
$ws.Rows.Item(257).Insert()

# Populate the new row 257 with the new weekly price-report record.
$ws.Cells.Item(257, 1).Value = 4
$ws.Cells.Item(257, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(257, 3).Value = "Los Lagos"
$ws.Cells.Item(257, 4).Value = 44782
$ws.Cells.Item(257, 5).Value = 10
$ws.Cells.Item(257, 6).Value = 100112037
$ws.Cells.Item(257, 7).Value = "Cebollín"
$ws.Cells.Item(257, 8).Value = "Sin especificar"
$ws.Cells.Item(257, 9).Value = "Primera"
$ws.Cells.Item(257, 10).Value = 180
$ws.Cells.Item(257, 11).Value = 9500
$ws.Cells.Item(257, 12).Value = 10000
$ws.Cells.Item(257, 13).Value = 9750
$ws.Cells.Item(257, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(257, 15).Value = "Región Metropolitana"
$ws.Cells.Item(257, 16).Value = 271
$ws.Cells.Item(257, 17).Value = 36
$ws.Cells.Item(257, 18).Value = "Hortaliza"
